# Update "Forecast Comparison" sheet: columns D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast) and G (Amazon P90
# Forecast) for rows 2-17, reflecting removal of the Auto Arima model
# from the forecast averaging (new forecast values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @{
    2  = @(18, 18, 29, 51)
    3  = @(19, 20, 32, 55)
    4  = @(25, 26, 42, 71)
    5  = @(27, 30, 46, 74)
    6  = @(29, 31, 48, 80)
    7  = @(28, 29, 46, 77)
    8  = @(29, 31, 48, 81)
    9  = @(29, 30, 48, 81)
    10 = @(29, 31, 48, 80)
    11 = @(28, 29, 46, 78)
    12 = @(28, 29, 46, 79)
    13 = @(30, 32, 50, 83)
    14 = @(29, 30, 47, 80)
    15 = @(29, 29, 47, 82)
    16 = @(29, 31, 49, 82)
    17 = @(28, 29, 46, 79)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
    $ws.Cells.Item($row, 6).Value = $values[2]
    $ws.Cells.Item($row, 7).Value = $values[3]
}
